$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; old D,E,F (arrecadado_sucesso, taxa_sucesso, media_sucesso)
# shift right to become E,F,G. This keeps the <cols> width metadata aligned with the new layout.
$ws.Range("D1").EntireColumn.Insert()

# --- Headers (row 1) ---
$ws.Range("D1").Value = "particip"
$ws.Range("E1").Value = "taxa_sucesso"
$ws.Range("F1").Value = "arrecadado_sucesso"
$ws.Range("G1").Value = "media_sucesso"
$ws.Range("H1").Value = "std_sucesso"
$ws.Range("I1").Value = "min_sucesso"
$ws.Range("J1").Value = "max_sucesso"
$ws.Range("K1").Value = "apoio_medio"
$ws.Range("L1").Value = "contribuicoes"
$ws.Range("M1").Value = "media_contribuicoes"

# --- Data rows ---
# row 2 - year 2016
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 15.15151515151515
$ws.Range("F2").Value = 426.4616061876675
$ws.Range("G2").Value = 85.29232123753351
$ws.Range("H2").Value = 85.51030885495558
$ws.Range("I2").Value = 7.154956142241136
$ws.Range("J2").Value = 226.3900843036052
$ws.Range("K2").Value = 18.5418089646812
$ws.Range("L2").Value = 23
$ws.Range("M2").Value = 4.6

# row 3 - year 2017
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 9.929078014184398
$ws.Range("F3").Value = 4736.523382339164
$ws.Range("G3").Value = 338.3230987385117
$ws.Range("H3").Value = 458.6831840070852
$ws.Range("I3").Value = 27.62335886703489
$ws.Range("J3").Value = 1809.09852121176
$ws.Range("K3").Value = 16.91615493692558
$ws.Range("L3").Value = 280
$ws.Range("M3").Value = 20

# row 4 - year 2018
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 14.66666666666667
$ws.Range("F4").Value = 7327.98576013602
$ws.Range("G4").Value = 666.180523648729
$ws.Range("H4").Value = 1084.782975461776
$ws.Range("I4").Value = 26.58043580770418
$ws.Range("J4").Value = 3475.049171548047
$ws.Range("K4").Value = 20.99709386858459
$ws.Range("L4").Value = 349
$ws.Range("M4").Value = 31.72727272727273

# row 5 - year 2019
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 18.82352941176471
$ws.Range("F5").Value = 3191.257392255826
$ws.Range("G5").Value = 199.4535870159891
$ws.Range("H5").Value = 234.097723686266
$ws.Range("I5").Value = 10.31772032536115
$ws.Range("J5").Value = 834.8528000913501
$ws.Range("K5").Value = 17.06554755217019
$ws.Range("L5").Value = 187
$ws.Range("M5").Value = 11.6875

# row 6 - year 2020
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 30.8411214953271
$ws.Range("F6").Value = 15024.86326411449
$ws.Range("G6").Value = 455.2988867913483
$ws.Range("H6").Value = 1126.99584249591
$ws.Range("I6").Value = 3.799754022893506
$ws.Range("J6").Value = 5087.076865717208
$ws.Range("K6").Value = 20.44199083552992
$ws.Range("L6").Value = 735
$ws.Range("M6").Value = 22.27272727272727

# row 7 - year 2021
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 21.42857142857143
$ws.Range("F7").Value = 2767.438420781004
$ws.Range("G7").Value = 131.7827819419526
$ws.Range("H7").Value = 138.1458157071481
$ws.Range("I7").Value = 5.763382152582333
$ws.Range("J7").Value = 538.4389998789497
$ws.Range("K7").Value = 24.70927161411611
$ws.Range("L7").Value = 112
$ws.Range("M7").Value = 5.333333333333333

# row 8 - year 2022
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 31.88405797101449
$ws.Range("F8").Value = 4495.625203875862
$ws.Range("G8").Value = 204.3466001761755
$ws.Range("H8").Value = 345.0624243778461
$ws.Range("I8").Value = 1.087396962410123
$ws.Range("J8").Value = 1594.029696524064
$ws.Range("K8").Value = 23.66118528355717
$ws.Range("L8").Value = 190
$ws.Range("M8").Value = 8.636363636363637

# row 9 - year 2023
$ws.Range("D9").Value = 100
$ws.Range("E9").Value = 39.47368421052632
$ws.Range("F9").Value = 5216.802725094768
$ws.Range("G9").Value = 173.8934241698256
$ws.Range("H9").Value = 195.9730930324503
$ws.Range("I9").Value = 2.022084306600051
$ws.Range("J9").Value = 657.0789958678034
$ws.Range("K9").Value = 15.71326122016496
$ws.Range("L9").Value = 332
$ws.Range("M9").Value = 11.06666666666667

# --- Number formats / styles ---
# D (particip) has no special formatting
$ws.Range("D2:D9").Style = "Normal"
# E (taxa_sucesso) keeps percent formatting
$ws.Range("E2:E9").NumberFormat = "0.00%"
# F (arrecadado_sucesso) and G (media_sucesso) use currency formatting
$ws.Range("F2:F9").NumberFormat = "R$ #,##0.00"
$ws.Range("G2:G9").NumberFormat = "R$ #,##0.00"
# H:M (std_sucesso .. media_contribuicoes) have no special formatting
$ws.Range("H2:M9").Style = "Normal"
